# OY-3070 validate presence of hakemus OID when henkilo OID present
#
# This adds new columns for person/application identification details
# (Kutsumanimi, Syntymäpaikka, Passin numero, Kansallinen ID-tunnus,
# Kaupunki ja maa, Hakemus-oid) to the Sheet1 export template, updates
# the sample applicant OID value and keeps the existing e-mail hyperlink
# but makes its display text explicit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting of the new trailing columns (AE:AJ) for all existing
# rows (1-10) by copying the look of the last currently formatted column
# (AD), so the new cells inherit identical header/body/blank styling.
$ws.Range("AD1:AD10").Copy()
$ws.Range("AE1:AJ10").PasteSpecial(-4122)

# New header row (row 1) labels for the newly appended columns.
$ws.Range("AE1").Value = "Kutsumanimi"
$ws.Range("AF1").Value = "Syntymäpaikka"
$ws.Range("AG1").Value = "Passin numero"
$ws.Range("AH1").Value = "Kansallinen ID-tunnus"
$ws.Range("AI1").Value = "Kaupunki ja maa"
$ws.Range("AJ1").Value = "Hakemus-oid"

# The sample applicant OID value used to be a plain placeholder; update it
# to look like a real OID-ish value.
$ws.Range("G2").Value = "Hakijaoid1"

# Populate the new sample-data row (row 2) values for the columns that
# should contain something ("Kutsumanimi" nickname and the "Hakemus-oid"
# application OID used by the OY-3070 validation). The remaining new
# columns (Syntymäpaikka, Passin numero, Kansallinen ID-tunnus, Kaupunki ja
# maa) stay blank, matching the other optional/unused sample fields.
$ws.Range("AE2").Value = "Tuomas"
$ws.Range("AJ2").Value = "Hakemus1"

# Keep the existing e-mail hyperlink on D2 but make the link's display
# text explicit (it previously relied on the cell text alone).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:tuomas.hakkarainen@example.com", "", "", "tuomas.hakkarainen@example.com")

Write-Host "Applied OY-3070 column additions and hakemus-oid/hakija-oid updates"
